# Insert two new product rows ("cabbiage" id=2 and "onion" id=19) above the
# existing "rice" row, pushing "rice"/"white rice"/"brown rice" down by two
# rows (old rows 3-5 -> new rows 5-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 3 (each Insert() pushes everything below down
# by one row, carrying along the date-column (F/I/L) number format).
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(4).Insert()

function Set-TextValue($range, $text) {
    # Force the cell to keep a numeric-looking string as literal text
    # (matching the workbook's existing "number stored as text" columns),
    # then drop back to the default style so no stray number format sticks.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Row 3: cabbiage ---------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "cabbiage"
Set-TextValue $ws.Range("C3") "1.00"
$ws.Range("D3").Value = "supplier_1"
Set-TextValue $ws.Range("E3") "50.00"
$ws.Range("F3").Value = 45069.645949074074
$ws.Range("J3").Value = "supplier_1"
Set-TextValue $ws.Range("K3") "50.00"
$ws.Range("L3").Value = 45069.645949074074

# --- Row 4: onion -------------------------------------------------------
$ws.Range("A4").Value = 19
$ws.Range("B4").Value = "onion"
Set-TextValue $ws.Range("C4") "1.00"
$ws.Range("D4").Value = "supplier_1"
Set-TextValue $ws.Range("E4") "50.00"
$ws.Range("F4").Value = 45069.645949074074
$ws.Range("J4").Value = "supplier_1"
Set-TextValue $ws.Range("K4") "50.00"
$ws.Range("L4").Value = 45069.645949074074
